$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","I","J","M","N")

$data = @{
    2 = @(2.907024114200851, 1.511183025706146, 0.5062788764427495, 0.152218325803819, 0.002634331936744694, 2.441756269691155, 0.04700329101950551, 1.175457296432327, 2.930400163586981)
    3 = @(2.786335481763786, 1.453549657153815, 0.5042188027763217, 0.1522659843660108, 0.002641978109459542, 2.412453084211762, 0.04720539459368123, 1.143435964053609, 2.927344641416909)
    4 = @(2.714299125305786, 1.419208695023485, 0.5031992914765908, 0.1523572754431157, 0.002646912661333376, 2.395631484059962, 0.04733921665385221, 1.124556620446512, 2.92618224882132)
    5 = @(2.685458671348101, 1.405474663956454, 0.5028453632690457, 0.152410051814126, 0.002648984065920356, 2.389068752388795, 0.04739620037303549, 1.117058275872054, 2.925886717108114)
    6 = @(2.680700698943724, 1.403209770388912, 0.5027903051549742, 0.1524197553029687, 0.002649331684118403, 2.387996601661342, 0.04740581056402249, 1.115824925185869, 2.925848366068138)
    7 = @(2.713908094125316, 1.419022423334525, 0.5031942693622256, 0.1523579241720583, 0.002646940351837071, 2.395541796461458, 0.0473399752312158, 1.124454706943837, 2.926177543565046)
    8 = @(2.864978712027948, 1.491092142976299, 0.5055175266139855, 0.1522218699154649, 0.002636918718904591, 2.431408228599153, 0.04707095905368153, 1.164253154238267, 2.929197688808785)
    9 = @(3.177867535666564, 1.640862032606833, 0.5120292843615175, 0.1524484762392859, 0.002619157723328402, 2.511135991058211, 0.04662046598712521, 1.248574750093653, 2.940844705335678)
    10 = @(3.418255672772716, 1.756254612747682, 0.5180196658576364, 0.1529177685499725, 0.002607246200960773, 2.575594640158258, 0.0463362555449276, 1.314463465635285, 2.952980064134863)
    11 = @(3.529979003534891, 1.809960280355199, 0.5210099662605217, 0.1531974900310971, 0.002602070988242611, 2.606230957416045, 0.04621707483944171, 1.345318165031983, 2.959297386329354)
    12 = @(3.572632826870176, 1.830475353494478, 0.5221807077741971, 0.153312975862967, 0.00260014601106633, 2.618023813429119, 0.04617339471005444, 1.357130761132069, 2.961805779015066)
    13 = @(3.563431051339592, 1.826049096954648, 0.5219268568205848, 0.1532876780796144, 0.002600559046957979, 2.615475449195472, 0.04618273751188262, 1.354580958583981, 2.961260360220876)
    14 = @(3.533481178161651, 1.811644482136444, 0.5211055132917579, 0.1532067992809836, 0.002601911923914336, 2.607197308595573, 0.04621345218238204, 1.346287406667898, 2.959501416205114)
    15 = @(3.515181333834789, 1.802844512341835, 0.5206074220301673, 0.1531585049655426, 0.002602745120138423, 2.602151732713452, 0.04623245469765536, 1.341224167843691, 2.958439186686888)
    16 = @(3.411002471289919, 1.752769524831933, 0.5178296002326874, 0.1529008240471548, 0.002607589290748168, 2.573619144431618, 0.04634424746740429, 1.312464932735423, 2.952583381092467)
    17 = @(3.347703151301232, 1.72236326637335, 0.5161935986881616, 0.1527597356818653, 0.002610623204198378, 2.556453662987138, 0.04641541580672914, 1.295049149199173, 2.949196242752606)
    18 = @(3.311517729765967, 1.704988351553254, 0.5152775653290007, 0.1526848177436833, 0.002612391155029386, 2.546704027926481, 0.04645730151424932, 1.285114918997053, 2.947322961613793)
    19 = @(3.299304035388673, 1.699124972357879, 0.5149716905109898, 0.152660521148114, 0.002612993697280185, 2.54342410032136, 0.04647164680481985, 1.281765537348804, 2.946701519527949)
    20 = @(3.354418384043015, 1.72558824508917, 0.5163651697714045, 0.1527741094956383, 0.002610297868192499, 2.558268157573906, 0.04640774135136816, 1.296894499560693, 2.949549045115276)
    21 = @(3.542268727003432, 1.815870610353727, 0.5213457182462378, 0.1532302955527882, 0.002601513609606778, 2.609623579335491, 0.04620439117819686, 1.348719919233687, 2.960014895505537)
    22 = @(3.667063010847414, 1.875913979997847, 0.5248246171092887, 0.1535841886819931, 0.002595975107011913, 2.644305066506945, 0.04607994663910198, 1.383341234185153, 2.967532751095916)
    23 = @(3.600270842767372, 1.843771499180207, 0.5229473027387428, 0.1533901957188988, 0.002598912657914434, 2.625691744631496, 0.04614559205680813, 1.364793910378594, 2.963457783865067)
    24 = @(3.351381785493629, 1.724129902156335, 0.5162875261105171, 0.1527675917991047, 0.00261044487877794, 2.557447454021144, 0.04641120795067621, 1.296059972800791, 2.949389312699196)
    25 = @(3.091404236406788, 1.599420114856059, 0.5100568904771734, 0.1523341484323133, 0.002623761665397017, 2.488545923991595, 0.04673410925480415, 1.225080218035771, 2.937071977213023)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

$wb.Save()